$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (cell E8) from "Good Morning" to "GIT UPDATE".
# This naturally causes the now-unused "Good Morning" shared string to be dropped and
# a new "GIT UPDATE" shared string to be appended, matching the target shared string table.
$ws.Range("E8").Value = "GIT UPDATE"

# Update the active selection shown in the sheet view to E8.
$ws.Range("E8").Select()
